# Auto-generated edit script
# Applies value changes per the commit diff to the Famfrit_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 999.5
$ws.Range("J19").Value = 999.5
$ws.Range("L19").Value = 999.5
$ws.Range("N19").Value = -1349.5
$ws.Range("H40").Value = 5953.846
$ws.Range("I40").Value = 3636.818
$ws.Range("K40").Value = 3636.818
$ws.Range("M40").Value = -3461.818
$ws.Range("H137").Value = 1187.5
$ws.Range("I137").Value = 1106.8182
$ws.Range("J137").Value = 1483.3334
$ws.Range("K137").Value = 3320.4546
$ws.Range("L137").Value = 4450.0002
$ws.Range("M137").Value = -770.4546
$ws.Range("N137").Value = -9550.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3015
$ws.Range("I2").Value = 2380
$ws.Range("J2").Value = 4158
$ws.Range("K2").Value = 2380
$ws.Range("L2").Value = 4158
$ws.Range("M2").Value = -2267
$ws.Range("N2").Value = -4384
$ws.Range("H74").Value = 37038750
$ws.Range("I74").Value = 45455628
$ws.Range("K74").Value = 45455628
$ws.Range("M74").Value = -45454754
$ws.Range("H77").Value = 37038750
$ws.Range("I77").Value = 45455628
$ws.Range("K77").Value = 227278140
$ws.Range("M77").Value = -227273772
$ws.Range("H88").Value = 15357
$ws.Range("I88").Value = 34585.332
$ws.Range("J88").Value = 3820
$ws.Range("K88").Value = 34585.332
$ws.Range("L88").Value = 3820
$ws.Range("M88").Value = -34179.332
$ws.Range("N88").Value = -4632
$ws.Range("H91").Value = 15357
$ws.Range("I91").Value = 34585.332
$ws.Range("J91").Value = 3820
$ws.Range("K91").Value = 34585.332
$ws.Range("L91").Value = 3820
$ws.Range("M91").Value = -33181.332
$ws.Range("N91").Value = -6628
$ws.Range("H102").Value = 291401.16
$ws.Range("I102").Value = 504452.5
$ws.Range("K102").Value = 504452.5
$ws.Range("M102").Value = -502830.5
$ws.Range("H116").Value = 3015
$ws.Range("I116").Value = 2380
$ws.Range("J116").Value = 4158
$ws.Range("K116").Value = 2380
$ws.Range("L116").Value = 4158
$ws.Range("M116").Value = -86
$ws.Range("N116").Value = -8746
$ws.Range("H132").Value = 27853936
$ws.Range("I132").Value = 12367.68
$ws.Range("K132").Value = 37103.04
$ws.Range("M132").Value = -34573.04

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3015
$ws.Range("I3").Value = 2380
$ws.Range("J3").Value = 4158
$ws.Range("K3").Value = 2380
$ws.Range("L3").Value = 4158
$ws.Range("M3").Value = -2266
$ws.Range("N3").Value = -4386
$ws.Range("H22").Value = 225
$ws.Range("I22").Value = 250.2
$ws.Range("J22").Value = 99
$ws.Range("K22").Value = 250.2
$ws.Range("L22").Value = 99
$ws.Range("M22").Value = -77.19999999999999
$ws.Range("N22").Value = -445
$ws.Range("H86").Value = 13432.75
$ws.Range("I86").Value = 23257.6
$ws.Range("J86").Value = 6415
$ws.Range("K86").Value = 23257.6
$ws.Range("L86").Value = 6415
$ws.Range("M86").Value = -22134.6
$ws.Range("N86").Value = -8661
$ws.Range("H89").Value = 13432.75
$ws.Range("I89").Value = 23257.6
$ws.Range("J89").Value = 6415
$ws.Range("K89").Value = 116288
$ws.Range("L89").Value = 32075
$ws.Range("M89").Value = -110672
$ws.Range("N89").Value = -43307
$ws.Range("H97").Value = 12851.272
$ws.Range("I97").Value = 12851.272
$ws.Range("K97").Value = 12851.272
$ws.Range("M97").Value = -11860.272
$ws.Range("H99").Value = 3680.8667
$ws.Range("I99").Value = 3337.7273
$ws.Range("K99").Value = 3337.7273
$ws.Range("M99").Value = -1839.7273
$ws.Range("H107").Value = 3380.2
$ws.Range("I107").Value = 2225.25
$ws.Range("J107").Value = 8000
$ws.Range("K107").Value = 2225.25
$ws.Range("L107").Value = 8000
$ws.Range("M107").Value = -305.25
$ws.Range("N107").Value = -11840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 8115.933
$ws.Range("I22").Value = 10274
$ws.Range("J22").Value = 3799.8
$ws.Range("K22").Value = 10274
$ws.Range("L22").Value = 3799.8
$ws.Range("M22").Value = -9924
$ws.Range("N22").Value = -4499.8
$ws.Range("H41").Value = 32841.715
$ws.Range("J41").Value = 32841.715
$ws.Range("L41").Value = 32841.715
$ws.Range("N41").Value = -33697.715
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250
$ws.Range("H107").Value = 358
$ws.Range("I107").Value = 221.16667
$ws.Range("K107").Value = 221.16667
$ws.Range("M107").Value = 1698.83333
$ws.Range("H131").Value = 25365.666
$ws.Range("J131").Value = 27438.8
$ws.Range("L131").Value = 27438.8
$ws.Range("N131").Value = -37518.8
$ws.Range("H141").Value = 109872.875
$ws.Range("J141").Value = 121283.29
$ws.Range("L141").Value = 121283.29
$ws.Range("N141").Value = -131643.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 77028
$ws.Range("I2").Value = 92
$ws.Range("J2").Value = 250134
$ws.Range("K2").Value = 552
$ws.Range("L2").Value = 1500804
$ws.Range("M2").Value = -439
$ws.Range("N2").Value = -1501030
$ws.Range("H12").Value = 1010.26666
$ws.Range("J12").Value = 967.63635
$ws.Range("L12").Value = 2902.90905
$ws.Range("N12").Value = -3248.90905
$ws.Range("H38").Value = 117.27273
$ws.Range("I38").Value = 18.5
$ws.Range("K38").Value = 55.5
$ws.Range("M38").Value = 291.5
$ws.Range("H80").Value = 502
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 502
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 26829
$ws.Range("J20").Value = 28804.572
$ws.Range("L20").Value = 28804.572
$ws.Range("N20").Value = -29294.572
$ws.Range("H21").Value = 12160.833
$ws.Range("I21").Value = 12160.833
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 12160.833
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -11987.833
$ws.Range("N21").ClearContents()
$ws.Range("H24").Value = 38812.832
$ws.Range("I24").Value = 21995
$ws.Range("J24").Value = 47221.75
$ws.Range("K24").Value = 21995
$ws.Range("L24").Value = 47221.75
$ws.Range("M24").Value = -21822
$ws.Range("N24").Value = -47567.75
$ws.Range("H30").Value = 12160.833
$ws.Range("I30").Value = 12160.833
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 12160.833
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -12055.833
$ws.Range("N30").ClearContents()
$ws.Range("H122").Value = 50001884
$ws.Range("I122").Value = 3253.5
$ws.Range("J122").Value = 62501540
$ws.Range("K122").Value = 9760.5
$ws.Range("L122").Value = 187504620
$ws.Range("M122").Value = -7310.5
$ws.Range("N122").Value = -187509520
$ws.Range("H132").Value = 9531.916999999999
$ws.Range("I132").Value = 9443.375
$ws.Range("J132").Value = 9709
$ws.Range("K132").Value = 28330.125
$ws.Range("L132").Value = 29127
$ws.Range("M132").Value = -25800.125
$ws.Range("N132").Value = -34187
$ws.Range("H133").Value = 195999
$ws.Range("J133").Value = 195999
$ws.Range("L133").Value = 195999
$ws.Range("N133").Value = -206119

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 928.6429000000001
$ws.Range("I22").Value = 1252.8889
$ws.Range("J22").Value = 345
$ws.Range("K22").Value = 1252.8889
$ws.Range("L22").Value = 345
$ws.Range("M22").Value = -957.8888999999999
$ws.Range("N22").Value = -935
$ws.Range("H26").Value = 26254
$ws.Range("I26").Value = 26254
$ws.Range("K26").Value = 26254
$ws.Range("M26").Value = -25959
$ws.Range("H27").Value = 928.6429000000001
$ws.Range("I27").Value = 1252.8889
$ws.Range("J27").Value = 345
$ws.Range("K27").Value = 1252.8889
$ws.Range("L27").Value = 345
$ws.Range("M27").Value = -1145.8889
$ws.Range("N27").Value = -559
$ws.Range("H68").Value = 2546.5
$ws.Range("I68").Value = 2495.5
$ws.Range("K68").Value = 2495.5
$ws.Range("M68").Value = -1746.5
$ws.Range("H71").Value = 2546.5
$ws.Range("I71").Value = 2495.5
$ws.Range("K71").Value = 12477.5
$ws.Range("M71").Value = -8733.5
$ws.Range("H100").Value = 3176.6155
$ws.Range("I100").Value = 2753.8333
$ws.Range("K100").Value = 2753.8333
$ws.Range("M100").Value = -2212.8333
$ws.Range("H122").Value = 3052184.2
$ws.Range("I122").Value = 3328.5925
$ws.Range("J122").Value = 8932120
$ws.Range("K122").Value = 9985.7775
$ws.Range("L122").Value = 26796360
$ws.Range("M122").Value = -7535.7775
$ws.Range("N122").Value = -26801260

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4465.2104
$ws.Range("I136").Value = 2987.6155
$ws.Range("J136").Value = 7666.6665
$ws.Range("K136").Value = 8962.8465
$ws.Range("L136").Value = 22999.9995
$ws.Range("M136").Value = -6412.8465
$ws.Range("N136").Value = -28099.9995
